# Update gh-pages output (合肥-漫展信息.xlsx) to the data generated at 456a3b4.
#
# Sheets (by index, 1-based, matching tab order):
#   1 = 展览      (Exhibitions)
#   2 = 演出      (Performances)
#   3 = 本地生活  (Local life) -- unchanged
#   4 = 全部类型  (All types, union of the other three sheets)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: apply the "想去人数" (F column) refresh that is common to every
# sheet that carries these particular exhibition rows (展览 + 全部类型 both
# list the same 21 pre-existing exhibition rows at the same row numbers).
# ---------------------------------------------------------------------------
function Update-ExhibitionCounts($ws) {
    $ws.Range("F2").Value = 235
    $ws.Range("F3").Value = 438
    $ws.Range("F4").Value = 13085
    $ws.Range("F5").Value = 1344
    $ws.Range("F6").Value = 220
    $ws.Range("F8").Value = 98
    $ws.Range("F10").Value = 231
    $ws.Range("F13").Value = 69
    $ws.Range("F16").Value = 47
    $ws.Range("F17").Value = 419
    $ws.Range("F18").Value = 5555
    $ws.Range("F19").Value = 110
    $ws.Range("F20").Value = 59
    $ws.Range("F21").Value = 962
}

# ---------------------------------------------------------------------------
# Helper: a brand new exhibition ("合肥·城市动漫节") was scraped in between
# the existing row 21 (Look Look动漫嘉年华, 2024-06-22) and what used to be
# row 22 (W·A首届童年怀旧only, 2024-07-20). Shift the old rows 22..lastRow
# down by one (copying A:I only, so no stray formatting/columns are
# introduced), then write the new row 22, then fix up the two rows whose
# "想去人数" counts were refreshed by the scraper on top of the shift, then
# restore column A's running index (which must stay equal to row# - 1, i.e.
# untouched by the content shift).
# ---------------------------------------------------------------------------
function Insert-CityAnimeFestival($ws, $lastRowBefore) {
    # Shift rows down, bottom-up so we never overwrite data we still need.
    for ($r = $lastRowBefore; $r -ge 22; $r--) {
        $src = $ws.Range("A" + $r + ":I" + $r)
        $dst = $ws.Range("A" + ($r + 1) + ":I" + ($r + 1))
        $src.Copy($dst)
    }

    # New row 22: 合肥·城市动漫节
    # B22 looks like an ISO date ("yyyy-mm-dd"); Excel would otherwise
    # silently reinterpret it as a date serial number on assignment (like
    # every other cell in this column, it must stay plain text). Force a
    # text number format for the assignment, then drop back to the sheet's
    # default formatting so the cell ends up indistinguishable from its
    # neighbours (matching B23/B24/... which carry no explicit style).
    $b22 = $ws.Range("B22")
    $b22.NumberFormat = "@"
    $b22.Value = "2024-06-22"
    $b22.ClearFormats()
    $ws.Range("C22").Value = "合肥·城市动漫节"
    $ws.Range("D22").Value = "包河经济开发区大连路与园博大道交口骆岗中央公园园博小镇一期S6区1号楼 大机库演艺中心"
    $ws.Range("E22").Value = "2024.06.22 10:00-06.23 16:30"
    $ws.Range("F22").Value = 6
    $ws.Range("G22").Value = 50
    $ws.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=85000"
    $ws.Range("I22").Value = "//i0.hdslb.com/bfs/openplatform/202404/r4N9T80s1713843349802.jpeg"

    # The scraper also bumped the "想去人数" of two of the shifted rows:
    #   安徽·MAX特摄only展 (now row 24): 130 -> 131
    #   合肥·第七届环形宇宙动漫游戏嘉年华 (now row 25): 147 -> 151
    $ws.Range("F24").Value = 131
    $ws.Range("F25").Value = 151

    # Column A is a plain running index (row# - 1); restore it for every
    # row that got shifted down (its value was copied along with A:I).
    for ($r = 23; $r -le ($lastRowBefore + 1); $r++) {
        $ws.Range("A" + $r).Value = $r - 1
    }
}

# ---------------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
Update-ExhibitionCounts $ws1
Insert-CityAnimeFestival $ws1 24

# ---------------------------------------------------------------------------
# Sheet 2: 演出 -- single row, its "想去人数" was refreshed too.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 13

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 -- union of all the other sheets; same exhibition rows
# as sheet 1 plus the 演出 row tacked on at the end (old row 25).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
Update-ExhibitionCounts $ws4
Insert-CityAnimeFestival $ws4 25

# Row 26 on sheet 4 (old row 25) is the 演出 sheet's single record; its
# "想去人数" was refreshed the same way as sheet 2's F2.
$ws4.Range("F26").Value = 13
